$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C2 value
$ws.Range("C2").Value = 15.01

# Add new ranking rows (67-97)
$ws.Range("A67").Value = "אביב ואסקז"
$ws.Range("B67").Value = 1
$ws.Range("A68").Value = "דן פימה"
$ws.Range("B68").Value = 1
$ws.Range("A69").Value = "עדן ורד מרי"
$ws.Range("B69").Value = 1
$ws.Range("A70").Value = "אביב ואסקז"
$ws.Range("B70").Value = 1
$ws.Range("A71").Value = "עדן ורד מרי"
$ws.Range("B71").Value = 1
$ws.Range("A72").Value = "רומי הרשקוביץ"
$ws.Range("B72").Value = 1
$ws.Range("A73").Value = "לידור אלשטיין"
$ws.Range("B73").Value = 1
$ws.Range("A74").Value = "אביב ואסקז"
$ws.Range("B74").Value = 1
$ws.Range("A75").Value = "אביב ואסקז"
$ws.Range("B75").Value = 1
$ws.Range("A76").Value = "לידור אלשטיין"
$ws.Range("B76").Value = 1
$ws.Range("A77").Value = "רומי הרשקוביץ"
$ws.Range("B77").Value = 1
$ws.Range("A78").Value = "לידור אלשטיין"
$ws.Range("B78").Value = 1
$ws.Range("A79").Value = "עדן ורד מרי"
$ws.Range("B79").Value = 1
$ws.Range("A80").Value = "רומי הרשקוביץ"
$ws.Range("B80").Value = 1
$ws.Range("A81").Value = "יהלי דוייב"
$ws.Range("B81").Value = 6
$ws.Range("A82").Value = "אביב ואסקז"
$ws.Range("B82").Value = 6
$ws.Range("A83").Value = "יער אלביר"
$ws.Range("B83").Value = 6
$ws.Range("A84").Value = "ירון גלפנד"
$ws.Range("B84").Value = 6
$ws.Range("A85").Value = "איתי הראל"
$ws.Range("B85").Value = 1
$ws.Range("A86").Value = "יהלי דוייב"
$ws.Range("B86").Value = 1
$ws.Range("A87").Value = "תאיו ורד"
$ws.Range("B87").Value = 1
$ws.Range("A88").Value = "יער אלביר"
$ws.Range("B88").Value = 1
$ws.Range("A89").Value = "עמית גורוביץ"
$ws.Range("B89").Value = 1
$ws.Range("A90").Value = "יולי יערי תליו"
$ws.Range("B90").Value = 1
$ws.Range("A91").Value = "מעיין סטרוזר"
$ws.Range("B91").Value = 1
$ws.Range("A92").Value = "ליהי בראל"
$ws.Range("B92").Value = 1
$ws.Range("A93").Value = "ירון גלפנד"
$ws.Range("B93").Value = 1
$ws.Range("A94").Value = "לינוי קוסטיקה"
$ws.Range("B94").Value = 1
$ws.Range("A95").Value = "תומר ששון"
$ws.Range("B95").Value = 1
$ws.Range("A96").Value = "אורי שטרנברג"
$ws.Range("B96").Value = 1
$ws.Range("A97").Value = "יובל סטרוזר"
$ws.Range("B97").Value = 1

# Update selection to C5 (view state)
$ws.Range("C5").Select() | Out-Null
